$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new values look like plain numbers (e.g. "304.90") must be
# forced to Text format first, otherwise Excel auto-converts them to floating
# point numbers and the trailing/formatting digits (and exact text) are lost.
# Values that already contain two "." separators (e.g. "22.447.87") are never
# parsed as numbers by Excel, so they do not need this treatment.
$textCells = @(
  "D5",
  "D6",
  "D7",
  "D8",
  "D9",
  "D10",
  "D11",
  "D12",
  "D13",
  "D14",
  "D15",
  "D16",
  "D18",
  "D19",
  "D20",
  "D21",
  "D22",
  "D23",
  "D25",
  "D26",
  "D27",
  "D28",
  "D31",
  "D32",
  "D33",
  "D34",
  "D35",
  "D36",
  "D37",
  "D39",
  "D40",
  "D41",
  "D42",
  "D43",
  "D44",
  "D45",
  "D46",
  "D47",
  "D48",
  "D49",
  "D51"
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '22.447.87'
$ws.Range("E2").Value = '  +9.05%  '
$ws.Range("D3").Value = '1.604.07'
$ws.Range("E3").Value = '  +8.34%  '
$ws.Range("E4").Value = '  -0.71%  '
$ws.Range("D5").Value = '304.90'
$ws.Range("E5").Value = '  +8.79%  '
$ws.Range("D6").Value = '0.9912'
$ws.Range("E6").Value = '  +2.04%  '
$ws.Range("D7").Value = '0.3694'
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("D8").Value = '0.3400'
$ws.Range("E8").Value = '  +10.25%  '
$ws.Range("D9").Value = '42.38'
$ws.Range("E9").Value = '  +5.59%  '
$ws.Range("D10").Value = '1.141'
$ws.Range("E10").Value = '  +7.48%  '
$ws.Range("D11").Value = '0.07078'
$ws.Range("E11").Value = '  +6.12%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '19.78'
$ws.Range("E13").Value = '  +9.06%  '
$ws.Range("D14").Value = '5.947'
$ws.Range("E14").Value = '  +7.46%  '
$ws.Range("D15").Value = '6.651'
$ws.Range("E15").Value = '  +7.00%  '
$ws.Range("D16").Value = '0.00001092'
$ws.Range("E16").Value = '  +5.93%  '
$ws.Range("D17").Value = '1.600.80'
$ws.Range("E17").Value = '  +7.90%  '
$ws.Range("D18").Value = '0.9914'
$ws.Range("E18").Value = '  +2.10%  '
$ws.Range("D19").Value = '0.06812'
$ws.Range("E19").Value = '  +14.38%  '
$ws.Range("D20").Value = '78.12'
$ws.Range("E20").Value = '  +11.75%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '16.15'
$ws.Range("E21").Value = '  +11.07%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.044'
$ws.Range("E22").Value = '  +9.76%  '
$ws.Range("D23").Value = '11.86'
$ws.Range("E23").Value = '  +7.12%  '
$ws.Range("D24").Value = '22.461.30'
$ws.Range("E24").Value = '  +8.88%  '
$ws.Range("D25").Value = '2.399'
$ws.Range("E25").Value = '  +5.53%  '
$ws.Range("D26").Value = '2.545'
$ws.Range("E26").Value = '  +20.09%  '
$ws.Range("D27").Value = '151.29'
$ws.Range("E27").Value = '  +6.45%  '
$ws.Range("D28").Value = '19.62'
$ws.Range("E28").Value = '  +13.28%  '
$ws.Range("D29").Value = '1.781.42'
$ws.Range("E29").Value = '  +8.23%  '
$ws.Range("E30").Value = '  +5.91%  '
$ws.Range("D31").Value = '4.181'
$ws.Range("E31").Value = '  +5.09%  '
$ws.Range("D32").Value = '6.129'
$ws.Range("E32").Value = '  +21.82%  '
$ws.Range("D33").Value = '0.9547'
$ws.Range("E33").Value = '  +15.60%  '
$ws.Range("D34").Value = '0.08270'
$ws.Range("E34").Value = '  +3.10%  '
$ws.Range("D35").Value = '1.636'
$ws.Range("E35").Value = '  +5.86%  '
$ws.Range("D36").Value = '5.286'
$ws.Range("E36").Value = '  +11.05%  '
$ws.Range("D37").Value = '12.01'
$ws.Range("E37").Value = '  +15.07%  '
$ws.Range("E38").Value = '  +3.59%  '
$ws.Range("D39").Value = '8.629'
$ws.Range("E39").Value = '  +12.50%  '
$ws.Range("D40").Value = '0.06109'
$ws.Range("E40").Value = '  +5.52%  '
$ws.Range("D41").Value = '0.02229'
$ws.Range("E41").Value = '  +8.81%  '
$ws.Range("D42").Value = '0.2030'
$ws.Range("E42").Value = '  +7.80%  '
$ws.Range("D43").Value = '0.9912'
$ws.Range("E43").Value = '  +2.12%  '
$ws.Range("D44").Value = '0.5935'
$ws.Range("E44").Value = '  +11.67%  '
$ws.Range("D45").Value = '3.844'
$ws.Range("E45").Value = '  +8.47%  '
$ws.Range("D46").Value = '13.18'
$ws.Range("E46").Value = '  +6.84%  '
$ws.Range("D47").Value = '0.5717'
$ws.Range("E47").Value = '  +9.70%  '
$ws.Range("D48").Value = '127.66'
$ws.Range("E48").Value = '  +7.25%  '
$ws.Range("D49").Value = '1.985'
$ws.Range("E49").Value = '  +8.54%  '
$ws.Range("E50").Value = '  +4.82%  '
$ws.Range("D51").Value = '73.97'
$ws.Range("E51").Value = '  +9.06%  '
